# fixbug tinh chiet khau don thu no
# Adds two new invoice rows (688, 689) into "Đơn sale chính" sheet,
# pushing the "Tổng" (total) row down and recomputing the totals,
# then updates the dependent summary values on the "Lương" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$ws2 = $wb.Worksheets.Item("Lương")

# Insert two new blank rows above the current "Tổng" row (row 6),
# which pushes it down to row 8.
$ws1.Range("A6:A7").EntireRow.Insert()

# The "Ngày thực hiện" column stores plain text dates like "08-27-2024"
# (not real Excel dates) - force text format so the assignment below
# doesn't get auto-coerced into a date serial number.
$ws1.Range("C6:C7").NumberFormat = "@"

# --- New row 6: HD-LUXURY 688 ---
$ws1.Cells.Item(6, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(6, 2).Value = 688
$ws1.Cells.Item(6, 3).Value = "08-27-2024"
$ws1.Cells.Item(6, 4).Value = "LONG XUYÊN"
$ws1.Cells.Item(6, 5).Value = "C.hạnh"
$ws1.Cells.Item(6, 6).Value = "Cá nhân"
$ws1.Cells.Item(6, 7).Value = "Cắt mí"
$ws1.Cells.Item(6, 8).Value = 5500000
$ws1.Cells.Item(6, 9).Value = ""
$ws1.Cells.Item(6, 10).Value = ""
$ws1.Cells.Item(6, 11).Value = 5500000
$ws1.Cells.Item(6, 12).Value = 5500000
$ws1.Cells.Item(6, 13).Value = 0.1
$ws1.Cells.Item(6, 14).Value = 550000

# --- New row 7: HD-LUXURY 689 ---
$ws1.Cells.Item(7, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(7, 2).Value = 689
$ws1.Cells.Item(7, 3).Value = "08-27-2024"
$ws1.Cells.Item(7, 4).Value = "LONG XUYÊN"
$ws1.Cells.Item(7, 5).Value = "Cầm dương"
$ws1.Cells.Item(7, 6).Value = "Cá nhân"
$ws1.Cells.Item(7, 7).Value = "Cắt mí"
$ws1.Cells.Item(7, 8).Value = 14000000
$ws1.Cells.Item(7, 9).Value = ""
$ws1.Cells.Item(7, 10).Value = ""
$ws1.Cells.Item(7, 11).Value = 14000000
$ws1.Cells.Item(7, 12).Value = 2000000
$ws1.Cells.Item(7, 13).Value = 0.15
$ws1.Cells.Item(7, 14).Value = 300000

# --- Row 8 ("Tổng") updated totals ---
$ws1.Cells.Item(8, 1).Value = "Tổng"
$ws1.Cells.Item(8, 2).Value = 6
$ws1.Cells.Item(8, 3).Value = ""
$ws1.Cells.Item(8, 4).Value = ""
$ws1.Cells.Item(8, 5).Value = ""
$ws1.Cells.Item(8, 6).Value = ""
$ws1.Cells.Item(8, 7).Value = ""
$ws1.Cells.Item(8, 8).Value = 39500000
$ws1.Cells.Item(8, 9).Value = ""
$ws1.Cells.Item(8, 10).Value = 3000000
$ws1.Cells.Item(8, 11).Value = 42500000
$ws1.Cells.Item(8, 12).Value = 26500000
$ws1.Cells.Item(8, 13).Value = 0
$ws1.Cells.Item(8, 14).Value = 2600000

# --- Update the dependent "Lương" sheet summary numbers ---
$ws2.Cells.Item(12, 2).Value = 24.5
$ws2.Cells.Item(13, 2).Value = 7000000.000000001
$ws2.Cells.Item(14, 2).Value = 2600000
$ws2.Cells.Item(32, 2).Value = 600000
$ws2.Cells.Item(34, 2).Value = 600000
